$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.824931474
$ws.Range("C2").Value = 6.544755634
$ws.Range("D2").Value = 1.479027748
$ws.Range("E2").Value = 0.5702867508
$ws.Range("F2").Value = 0.7198241599999999
$ws.Range("G2").Value = 0.9530009179505683
$ws.Range("H2").Value = 0.1516748068629467
$ws.Range("I2").Value = 8.728778021537353

$ws.Range("B3").Value = 6.544755634
$ws.Range("C3").Value = 7.264579827
$ws.Range("D3").Value = 0.5702867508
$ws.Range("E3").Value = 0.3810100555
$ws.Range("F3").Value = 0.719824193
$ws.Range("G3").Value = 0.4033135387425803
$ws.Range("H3").Value = 0.06418934330676629
$ws.Range("I3").Value = 8.72877762137065

$ws.Range("B4").Value = 7.804447965
$ws.Range("C4").Value = 8.524272157
$ws.Range("D4").Value = 0.1952342987
$ws.Range("E4").Value = 0.1958169937
$ws.Range("F4").Value = 0.7198241920000008
$ws.Range("G4").Value = -0.002980148346680622
$ws.Range("H4").Value = -0.0004743053405213604
$ws.Range("I4").Value = 8.728777633496902
